$wb = $excel.ActiveWorkbook

# FAVÖK row values (columns B..BV), in order
$values = @(0,0,0,0,34201000,34074000,0,36471000,52782000,49337000,0,58345000,77342000,69501000,61872000,102240000,115114000,96842000,0,-1935000,27201000,27263000,0,43873000,69143000,68572000,0,109977000,142577000,124764000,0,112623000,90416000,89688000,153934000,64251000,84501000,120446000,245303000,271991000,146358000,227865000,0,465619000,481949000,501829000,0,482645000,731915000,1006943000,0,986773000,1384589000,1343452000,0,1995597000,1886064000,1825506000,0,2046228000,1661788000,1618206000,0,1804652000,1916568000,2007054000,0,2719768000,2684064000,2990391000,3427066000,3854707000,3335573000)

$targetSheets = @("DATA_RAW", "gelir tablosu (çeyreklik)")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find first empty row right below the current used range in column A
    $newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

    $ws.Cells.Item($newRow, 1).Value = "FAVÖK"

    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($newRow, $i + 2).Value = $values[$i]
    }
}
